$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.751.77"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "3.146.59"
$ws.Range("E3").Value = "  -7.96%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "564.58"
$ws.Range("E5").Value = "  -3.36%  "
$ws.Range("D6").Value = "170.77"
$ws.Range("E6").Value = "  -4.58%  "
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "3.144.16"
$ws.Range("E9").Value = "  -7.85%  "
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -5.77%  "
$ws.Range("D11").Value = "6.56"
$ws.Range("E11").Value = "  -5.46%  "
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -4.73%  "
$ws.Range("D13").Value = "3.694.18"
$ws.Range("E13").Value = "  -7.88%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "27.19"
$ws.Range("E15").Value = "  -7.89%  "
$ws.Range("D16").Value = "64.601.37"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("E17").Value = "  -5.76%  "
$ws.Range("D18").Value = "3.155.93"
$ws.Range("E18").Value = "  -7.71%  "
$ws.Range("D19").Value = "5.72"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").Value = "12.91"
$ws.Range("E20").Value = "  -6.41%  "
$ws.Range("D21").Value = "355.82"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "7.22"
$ws.Range("E22").Value = "  -5.05%  "
$ws.Range("D24").Value = "68.64"
$ws.Range("E24").Value = "  -5.72%  "
$ws.Range("D25").Value = "0.499"
$ws.Range("E25").Value = "  -6.90%  "
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -7.62%  "
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "1.89"
$ws.Range("E31").Value = "  -5.03%  "
$ws.Range("D32").Value = "5.37"
$ws.Range("E32").Value = "  -7.81%  "
$ws.Range("D33").Value = "21.97"
$ws.Range("E33").Value = "  -6.33%  "
$ws.Range("D34").Value = "6.64"
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  -7.10%  "
$ws.Range("D37").Value = "153.76"
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("D38").Value = "0.829"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").Value = "26.03"
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("D40").Value = "1.73"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("D41").Value = "2.52"
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("D42").Value = "2.653.29"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").Value = "4.18"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("D44").Value = "6.04"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").Value = "24.23"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "39.17"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0654"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("D48").Value = "321.48"
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D49").Value = "0.0271"
$ws.Range("E49").Value = "  -5.08%  "
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.11%  "
